$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 19, pushing rows 19:25 down to 20:26
$ws.Rows.Item(19).Insert()

# Populate the new row 19 with the latest weekly record
$ws.Cells.Item(19, 1).Value = 11
$ws.Cells.Item(19, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(19, 3).Value = "Bíobío"
$ws.Cells.Item(19, 4).Value = 44523
$ws.Cells.Item(19, 4).NumberFormat = $ws.Cells.Item(20, 4).NumberFormat
$ws.Cells.Item(19, 5).Value = 8
$ws.Cells.Item(19, 6).Value = 100114007
$ws.Cells.Item(19, 7).Value = "Jengibre"
$ws.Cells.Item(19, 8).Value = "Sin especificar"
$ws.Cells.Item(19, 9).Value = "Primera"
$ws.Cells.Item(19, 10).Value = 40
$ws.Cells.Item(19, 11).Value = 15000
$ws.Cells.Item(19, 12).Value = 16000
$ws.Cells.Item(19, 13).Value = 15500
$ws.Cells.Item(19, 14).Value = "$/caja 13 kilos"
$ws.Cells.Item(19, 15).Value = "Perú"
$ws.Cells.Item(19, 16).Value = 1192
$ws.Cells.Item(19, 17).Value = 13
$ws.Cells.Item(19, 18).Value = "Hortaliza"
